# Update the "Correspond Handoff Datetime" (E) and "Correspond Handback
# DateTime" (H) timestamps for the be0bb493-... and d0b3b94b-... rows
# (which were handed off/back together and so share identical timestamps)
# on both the zh-cn and de-de sheets, to reflect the newly generated
# handback report times.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-24 04:22:08"
$wsZhCn.Range("H3").Value = "2016-03-24 04:22:33"
$wsZhCn.Range("E4").Value = "2016-03-24 04:22:08"
$wsZhCn.Range("H4").Value = "2016-03-24 04:22:33"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-24 04:22:12"
$wsDeDe.Range("H3").Value = "2016-03-24 04:22:40"
$wsDeDe.Range("E4").Value = "2016-03-24 04:22:12"
$wsDeDe.Range("H4").Value = "2016-03-24 04:22:40"
